$d = $word.ActiveDocument

# --- Simple single-value text replacements (unique across the document) ---
$d.Content.Find.Execute("99.97", $true, $false, $false, $false, $false, $true, 1, $false, "0M", 2) | Out-Null
$d.Content.Find.Execute("1.55", $true, $false, $false, $false, $false, $true, 1, $false, "0M", 2) | Out-Null
$d.Content.Find.Execute("4514", $true, $false, $false, $false, $false, $true, 1, $false, "0M", 2) | Out-Null
$d.Content.Find.Execute("1555", $true, $false, $false, $false, $false, $true, 1, $false, "1608", 2) | Out-Null
$d.Content.Find.Execute("0.01772", $true, $false, $false, $false, $false, $true, 1, $false, "0.01763", 2) | Out-Null
$d.Content.Find.Execute("0.00611", $true, $false, $false, $false, $false, $true, 1, $false, "0.00624", 2) | Out-Null
$d.Content.Find.Execute("1.13972", $true, $false, $false, $false, $false, $true, 1, $false, "1.54813", 2) | Out-Null

# --- Collapse the last three multi-column (tab-separated) rows down to a single value each ---
$t = $d.Tables.Item(1)
$rowCount = $t.Rows.Count
$t.Cell($rowCount - 2, 1).Range.Text = "99.97"
$t.Cell($rowCount - 1, 1).Range.Text = "1.55"
$t.Cell($rowCount, 1).Range.Text = "4514"
